# Updated cryptos list on Thu Jul 11 06:57:28 UTC 2024 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values
# for the crypto ranking table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.024.69'
$ws.Range("E2").Value = '  -2.07%  '
$ws.Range("D3").Value = '3.108.08'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.60'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.70'
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.108.19'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.445'
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("E10").Value = '  -3.39%  '
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").Value = '3.638.99'
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("E14").Value = '  +3.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.64'
$ws.Range("E15").Value = '  -5.88%  '
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").Value = '58.062.17'
$ws.Range("E17").Value = '  -1.90%  '
$ws.Range("D18").Value = '3.078.86'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.73'
$ws.Range("E20").Value = '  -3.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.98'
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '343.17'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.514'
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("E25").Value = '  +2.24%  '
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '0.0₃0926'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.37'
$ws.Range("E30").Value = '  -7.35%  '
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.01'
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("E34").Value = '  -3.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.88'
$ws.Range("E35").Value = '  +2.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.64'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.17'
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.21'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("E39").Value = '  -4.69%  '
$ws.Range("E40").Value = '  -3.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.60'
$ws.Range("E41").Value = '  +8.14%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.687'
$ws.Range("E43").Value = '  +2.68%  '
$ws.Range("D44").Value = '3.146.87'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").Value = '2.269.78'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.994'
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.13'
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.52'
$ws.Range("E51").Value = '  -2.23%  '
